$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q3" sheet, duplicating the format of the existing
#    "2022-Q2" sheet (so fonts/borders/column layout/styles match exactly),
#    then trim it down to the 5 rows needed and fill in the 2022-Q3 figures.
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy() with a single "Before" argument inserts the duplicate right before
# the target sheet, pushing "2022-Q2" (and everything after it) one slot to
# the right.
$q2Sheet.Copy($q2Sheet)

$q3Sheet = $wb.Worksheets.Item("2022-Q2 (2)")
$q3Sheet.Name = "2022-Q3"

# The duplicated sheet still has all 20 rows (1 header + 19 funds) from
# 2022-Q2; drop everything below row 5 so only 4 fund rows remain.
$q3Sheet.Range("A6:H20").EntireRow.Delete()

# Columns B (fund code, e.g. "005402") and D:G (decimal-looking percentages)
# must stay text so leading zeros / trailing zeros survive - force text
# format before assigning the values (NumberFormat needs to be set per
# contiguous block; comma-separated multi-area ranges don't apply reliably).
$q3Sheet.Range("B2:B5").NumberFormat = "@"
$q3Sheet.Range("D2:G5").NumberFormat = "@"

# Fill in the 2022-Q3 fund data (overwriting the copied 2022-Q2 values).
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "005402"
$q3Sheet.Range("C2").Value = "广发资源优选股票A"
$q3Sheet.Range("D2").Value = "7.40"
$q3Sheet.Range("E2").Value = "92.95"
$q3Sheet.Range("F2").Value = "5.39"
$q3Sheet.Range("G2").Value = "0.3989"
$q3Sheet.Range("H2").Value = 9

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "013422"
$q3Sheet.Range("C3").Value = "太平智行三个月定期开放混合"
$q3Sheet.Range("D3").Value = "6.38"
$q3Sheet.Range("E3").Value = "83.22"
$q3Sheet.Range("F3").Value = "5.29"
$q3Sheet.Range("G3").Value = "0.3375"
$q3Sheet.Range("H3").Value = 4

$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").Value = "009794"
$q3Sheet.Range("C4").Value = "太平智选一年定期开放股票"
$q3Sheet.Range("D4").Value = "5.02"
$q3Sheet.Range("E4").Value = "85.11"
$q3Sheet.Range("F4").Value = "4.78"
$q3Sheet.Range("G4").Value = "0.2400"
$q3Sheet.Range("H4").Value = 5

$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").Value = "010235"
$q3Sheet.Range("C5").Value = "广发资源优选股票C"
$q3Sheet.Range("D5").Value = "2.27"
$q3Sheet.Range("E5").Value = "92.95"
$q3Sheet.Range("F5").Value = "5.39"
$q3Sheet.Range("G5").Value = "0.1224"
$q3Sheet.Range("H5").Value = 9

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 above
#    the 2022-Q2 row and renumber the existing index column.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# The inserted row copies formatting from the header row above (bold,
# centered, bordered) for columns B:D - clear that so it matches the plain
# (unstyled) look of the other data rows.
$summary.Range("B2:D2").ClearFormats()

# Column A keeps the bordered/centered "index" style used by the other rows
# in this column; copy it from the row just below (now row 3).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 1.1

# Renumber the index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
